# Project DesignFirst save: update the "From" value for rule R30 (row 10)
# in the rules table from 18 to 100.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("C10").Value = 100
